$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Trajectory Ordering" section header (row 34), matching the
# underlined section-header look used by "Debug" (A31) etc.
$ws.Range("A34").Value = "Trajectory Ordering"
$ws.Range("A34").Font.Underline = 2

# Add its "Contour First?" parameter row (row 35), same look as the other
# parameter rows in column A (bold-free, default column font/style).
$ws.Range("A35").Value = "Contour First?"
$ws.Range("B35").Value = "Yes"

# Scroll the sheet view down so the new rows are visible (best effort —
# mirrors the author scrolling to topLeftCell A17 before saving).
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 17
